# Adds a bit of i18n: label the "Learning Factor" sheet's values with
# Russian-language descriptions in column A, pushing the existing values
# into column B, and add a third row for the output-layer dimensionality.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Learning Factor")

# Preserve the existing numeric values before we shuffle columns.
$learningRate = $ws.Range("A1").Value2
$inputCount   = $ws.Range("A2").Value2
$outputDim    = $ws.Range("B2").Value2

# Column B now holds the values; column A gets the new descriptive labels.
$ws.Range("B1").Value = $learningRate
$ws.Range("A1").Value = "Коэффициент скорости обучения"

$ws.Range("B2").Value = $inputCount
$ws.Range("A2").Value = "Количество входов нейронной сети"

$ws.Range("A3").Value = "Размерность выходного слоя"
$ws.Range("B3").Value = $outputDim

$ws.Range("N8").Select()
